# Apply the "Updated cryptos list" GitHub Actions refresh:
# updates Price (D) and Volume(1h) (E) columns for each coin row,
# and swaps the EnergySwap / Decentraland rows (47 and 48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.183.89"
$ws.Range("E2").Value = "  -5.49%  "
$ws.Range("D3").Value = "1.835.83"
$ws.Range("E3").Value = "  -5.13%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "'330.15"
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "'0.4600"
$ws.Range("E7").Value = "  -4.84%  "
$ws.Range("D8").Value = "'0.3857"
$ws.Range("E8").Value = "  -6.10%  "
$ws.Range("D9").Value = "'46.07"
$ws.Range("E9").Value = "  -3.87%  "
$ws.Range("D10").Value = "'0.07852"
$ws.Range("E10").Value = "  -3.82%  "
$ws.Range("D11").Value = "'0.9582"
$ws.Range("E11").Value = "  -5.41%  "
$ws.Range("D12").Value = "'21.86"
$ws.Range("E12").Value = "  -7.63%  "
$ws.Range("D13").Value = "1.836.58"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").Value = "'5.707"
$ws.Range("E14").Value = "  -6.08%  "
$ws.Range("D15").Value = "'6.905"
$ws.Range("E15").Value = "  -5.02%  "
$ws.Range("D16").Value = "'0.06852"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "'86.85"
$ws.Range("E18").Value = "  -4.39%  "
$ws.Range("D19").Value = "'0.000009928"
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("D20").Value = "'16.89"
$ws.Range("E20").Value = "  -4.95%  "
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").Value = "28.216.93"
$ws.Range("E22").Value = "  -5.34%  "
$ws.Range("D23").Value = "'5.335"
$ws.Range("E23").Value = "  -5.16%  "
$ws.Range("D24").Value = "'10.97"
$ws.Range("E24").Value = "  -7.58%  "
$ws.Range("D25").Value = "'2.132"
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("D26").Value = "2.054.59"
$ws.Range("E26").Value = "  -4.73%  "
$ws.Range("D27").Value = "'153.08"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "'19.18"
$ws.Range("E28").Value = "  -4.41%  "
$ws.Range("D29").Value = "'5.705"
$ws.Range("E29").Value = "  -13.23%  "
$ws.Range("D30").Value = "'1.978"
$ws.Range("E30").Value = "  -5.39%  "
$ws.Range("D31").Value = "'116.83"
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").Value = "'0.9406"
$ws.Range("E32").Value = "  -6.54%  "
$ws.Range("D33").Value = "'0.09272"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D34").Value = "'5.272"
$ws.Range("E34").Value = "  -4.95%  "
$ws.Range("D35").Value = "'3.445"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("D36").Value = "'1.322"
$ws.Range("E36").Value = "  -6.28%  "
$ws.Range("D37").Value = "'0.05995"
$ws.Range("E37").Value = "  -8.69%  "
$ws.Range("D38").Value = "'0.02150"
$ws.Range("E38").Value = "  -5.84%  "
$ws.Range("E39").Value = "  -4.72%  "
$ws.Range("D40").Value = "'1.001"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").Value = "'7.607"
$ws.Range("E41").Value = "  -4.34%  "
$ws.Range("D42").Value = "'0.5612"
$ws.Range("E42").Value = "  -6.17%  "
$ws.Range("D43").Value = "'9.969"
$ws.Range("E43").Value = "  -7.10%  "
$ws.Range("D44").Value = "'0.1773"
$ws.Range("E44").Value = "  -4.02%  "
$ws.Range("D45").Value = "'1.247"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").Value = "'2.248"
$ws.Range("E46").Value = "  -9.20%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'11.58"
$ws.Range("E47").Value = "  -6.40%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5273"
$ws.Range("E48").Value = "  -5.09%  "
$ws.Range("D49").Value = "'0.07019"
$ws.Range("E49").Value = "  -6.14%  "
$ws.Range("D50").Value = "'1.827"
$ws.Range("E50").Value = "  -7.90%  "
$ws.Range("D51").Value = "'112.63"
$ws.Range("E51").Value = "  -3.60%  "
